$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: paragraph "<head><m><pa>Skirret</pa> root</m></head>"
#   <m><pa>  ->  <pa>
#   remove the "</m>" run entirely
# ---------------------------------------------------------------------------

$anchor1 = $d.Content
$anchor1.Find.Execute("Skirret", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para1 = $anchor1.Paragraphs(1).Range
$p1Start = $para1.Start
$p1End = $para1.End

# <m><pa> -> <pa>
$scope1a = $d.Range($p1Start, $p1End)
$scope1a.Find.Execute("<m><pa>", $true, $false, $false, $false, $false, $true, 1, $false, "<pa>", 2)

# the paragraph shrank by 2 characters ("<m>" removed, 3 chars) -- recompute end
$p1End = $p1End - 3

# remove the "</m>" run
$scope1b = $d.Range($p1Start, $p1End)
$scope1b.Find.Execute("</m>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$scope1b.Delete()

# ---------------------------------------------------------------------------
# Part 2: paragraph about planting ("<ab>They want to be planted ... goodness.")
#   "They want to be planted in a very humid place where with such"
#       -> "They want to be planted " + <env> + "in a very humid place where with such"
#   <fr>  -> <oc>
#   </fr> -> </oc>
#   " fountain one can water it often, for by this method they are tender. ..."
#       -> " fountain one can water it often" + </env> + ", for by this method ..."
# ---------------------------------------------------------------------------

$anchor2 = $d.Content
$anchor2.Find.Execute("fountain", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para2 = $anchor2.Paragraphs(1).Range
$p2Start = $para2.Start
$p2End = $para2.End

# <fr> -> <oc>
$scope2a = $d.Range($p2Start, $p2End)
$scope2a.Find.Execute("<fr>", $true, $false, $false, $false, $false, $true, 1, $false, "<oc>", 2)

# </fr> -> </oc>
$scope2b = $d.Range($p2Start, $p2End)
$scope2b.Find.Execute("</fr>", $true, $false, $false, $false, $false, $true, 1, $false, "</oc>", 2)

# locate the (now) "<oc>" tag run -- we'll copy its run formatting (Courier New / blue / sz18)
$tagScope = $d.Range($p2Start, $p2End)
$tagScope.Find.Execute("<oc>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tagScope.Copy()

# split point #1: right after "They want to be planted "
$split1 = $d.Range($p2Start, $p2End)
$split1.Find.Execute("They want to be planted ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$split1.Collapse(0)
$split1Pos = $split1.Start
$split1.Paste()
$envOpenLen = 4   # length of "<oc>" text that got pasted
$newRun1 = $d.Range($split1Pos, $split1Pos + $envOpenLen)
$newRun1.Text = "<env>"

# paragraph grew by (len("<env>") - len("<oc>")) = 5 - 4 = 1
$p2End = $p2End + 1

# split point #2: right after " fountain one can water it often"
$split2 = $d.Range($p2Start, $p2End)
$split2.Find.Execute(" fountain one can water it often", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$split2.Collapse(0)
$split2Pos = $split2.Start
$split2.Paste()
$newRun2 = $d.Range($split2Pos, $split2Pos + $envOpenLen)
$newRun2.Text = "</env>"
